# sample_barang.xlsx now sources kode_barang as a generated UUID at
# runtime instead of shipping static sample codes, so the column's
# sample values are cleared out.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B4").ClearContents()

# Column A and B now share the same (narrower) width.
$ws.Range("B:B").ColumnWidth = 11.83

# Reflect the author's last selection before saving.
$ws.Range("J3").Select()
